$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Set header values for the new row
$ws.Range("A1").Value = "Question"
$ws.Range("B1").Value = "Answer"

# Update the selection to C6 as in the target workbook
$ws.Range("C6").Select()
